# Add a new "2021" column (O) to the table, mirroring the existing
# "2020" column (N) formatting, and update the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell O4: 2021, formatted like N4 (2020) ---
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O4").Value = 2021

# --- Data cell O5: 515, formatted like N5 (534) ---
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O5").Value = 515

# --- Update the view's selection/active cell (also drops topLeftCell) ---
$ws.Range("P12").Select() | Out-Null
